$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both have the "想去人数" (F column) figures updated.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1493
    $ws.Range("F3").Value = 3150
    $ws.Range("F5").Value = 895
}
